$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking values are written as exact text (avoid Excel auto-converting
# strings like "324.52", "-2.55%", "6" into numeric/percentage values).
$numericTextCells = @("D2", "E2", "G2", "D3", "E3", "G3", "D4", "E4", "G4", "D5", "E5", "G5", "D6", "E6", "G6", "D7", "E7", "G7", "D8", "E8", "G8", "D9", "E9", "G9", "E10", "G10", "D11", "E11", "G11", "D12", "E12", "G12", "D13", "E13", "G13", "D14", "E14", "G14", "D15", "E15", "G15", "E16", "G16", "D17", "E17", "G17", "D18", "E18", "G18", "D19", "E19", "G19", "D20", "E20", "G20", "D21", "E21", "G21", "D22", "E22", "G22", "D23", "E23", "G23", "D24", "E24", "G24", "D25", "E25", "G25", "D26", "E26", "G26", "D27", "E27", "G27", "G28", "G29", "G30", "G31", "G32", "G33", "G34", "G35", "G36", "G37", "G38", "D39", "E39", "G39", "D40", "E40", "G40", "D41", "E41", "G41", "D42", "E42", "G42", "D43", "E43", "G43", "D44", "E44", "G44", "D45", "E45", "G45", "D46", "E46", "G46", "E47", "G47", "D48", "E48", "G48", "D49", "E49", "G49", "D50", "E50", "G50", "D51", "E51", "G51")
foreach ($addr in $numericTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Updated price / volume / hora values
$ws.Range("D2").Value = "324.52"
$ws.Range("E2").Value = "-2.55%"
$ws.Range("G2").Value = "6"
$ws.Range("D3").Value = "44.66"
$ws.Range("E3").Value = "1.68%"
$ws.Range("G3").Value = "6"
$ws.Range("D4").Value = "5.499"
$ws.Range("E4").Value = "-4.66%"
$ws.Range("G4").Value = "6"
$ws.Range("D5").Value = "0.08035"
$ws.Range("E5").Value = "-3.73%"
$ws.Range("G5").Value = "6"
$ws.Range("D6").Value = "8.649"
$ws.Range("E6").Value = "-1.91%"
$ws.Range("G6").Value = "6"
$ws.Range("D7").Value = "4.305"
$ws.Range("E7").Value = "-4.40%"
$ws.Range("G7").Value = "6"
$ws.Range("D8").Value = "1.894"
$ws.Range("E8").Value = "-4.14%"
$ws.Range("G8").Value = "6"
$ws.Range("D9").Value = "2.713"
$ws.Range("E9").Value = "-6.46%"
$ws.Range("G9").Value = "6"
$ws.Range("E10").Value = "0.61%"
$ws.Range("G10").Value = "6"
$ws.Range("D11").Value = "0.1178"
$ws.Range("E11").Value = "-5.17%"
$ws.Range("G11").Value = "6"
$ws.Range("D12").Value = "0.1879"
$ws.Range("E12").Value = "-4.00%"
$ws.Range("G12").Value = "6"
$ws.Range("D13").Value = "0.09949"
$ws.Range("E13").Value = "3.21%"
$ws.Range("G13").Value = "6"
$ws.Range("D14").Value = "0.04189"
$ws.Range("E14").Value = "6.32%"
$ws.Range("G14").Value = "6"
$ws.Range("D15").Value = "0.1066"
$ws.Range("E15").Value = "-0.20%"
$ws.Range("G15").Value = "6"
$ws.Range("E16").Value = "-2.48%"
$ws.Range("G16").Value = "6"
$ws.Range("D17").Value = "0.005875"
$ws.Range("E17").Value = "-2.34%"
$ws.Range("G17").Value = "6"
$ws.Range("D18").Value = "0.004451"
$ws.Range("E18").Value = "1.48%"
$ws.Range("G18").Value = "6"
$ws.Range("D19").Value = "3.598"
$ws.Range("E19").Value = "2.67%"
$ws.Range("G19").Value = "6"
$ws.Range("D20").Value = "0.3486"
$ws.Range("E20").Value = "-0.66%"
$ws.Range("G20").Value = "6"
$ws.Range("D21").Value = "8.454"
$ws.Range("E21").Value = "-6.43%"
$ws.Range("G21").Value = "6"
$ws.Range("D22").Value = "0.1373"
$ws.Range("E22").Value = "0.16%"
$ws.Range("G22").Value = "6"
$ws.Range("D23").Value = "0.2643"
$ws.Range("E23").Value = "2.78%"
$ws.Range("G23").Value = "6"
$ws.Range("D24").Value = "0.04253"
$ws.Range("E24").Value = "-3.60%"
$ws.Range("G24").Value = "6"
$ws.Range("D25").Value = "0.001240"
$ws.Range("E25").Value = "-1.45%"
$ws.Range("G25").Value = "6"
$ws.Range("D26").Value = "0.0001234"
$ws.Range("E26").Value = "3.63%"
$ws.Range("G26").Value = "6"
$ws.Range("D27").Value = "0.0004003"
$ws.Range("E27").Value = "0.25%"
$ws.Range("G27").Value = "6"
$ws.Range("G28").Value = "6"
$ws.Range("G29").Value = "6"
$ws.Range("G30").Value = "6"
$ws.Range("G31").Value = "6"
$ws.Range("G32").Value = "6"
$ws.Range("G33").Value = "6"
$ws.Range("G34").Value = "6"
$ws.Range("G35").Value = "6"
$ws.Range("G36").Value = "6"
$ws.Range("G37").Value = "6"
$ws.Range("G38").Value = "6"
$ws.Range("D39").Value = "0.02629"
$ws.Range("E39").Value = "-7.48%"
$ws.Range("G39").Value = "6"
$ws.Range("D40").Value = "0.05477"
$ws.Range("E40").Value = "-4.59%"
$ws.Range("G40").Value = "6"
$ws.Range("D41").Value = "0.007681"
$ws.Range("E41").Value = "-3.02%"
$ws.Range("G41").Value = "6"
$ws.Range("D42").Value = "0.1392"
$ws.Range("E42").Value = "-2.46%"
$ws.Range("G42").Value = "6"
$ws.Range("D43").Value = "0.007212"
$ws.Range("E43").Value = "-20.49%"
$ws.Range("G43").Value = "6"
$ws.Range("D44").Value = "0.002054"
$ws.Range("E44").Value = "-2.31%"
$ws.Range("G44").Value = "6"
$ws.Range("D45").Value = "0.009190"
$ws.Range("E45").Value = "-9.78%"
$ws.Range("G45").Value = "6"
$ws.Range("D46").Value = "0.00007122"
$ws.Range("E46").Value = "-1.78%"
$ws.Range("G46").Value = "6"
$ws.Range("E47").Value = "0.26%"
$ws.Range("G47").Value = "6"
$ws.Range("D48").Value = "0.003488"
$ws.Range("E48").Value = "9.15%"
$ws.Range("G48").Value = "6"
$ws.Range("D49").Value = "0.002278"
$ws.Range("E49").Value = "-0.09%"
$ws.Range("G49").Value = "6"
$ws.Range("D50").Value = "0.00002108"
$ws.Range("E50").Value = "0.26%"
$ws.Range("G50").Value = "6"
$ws.Range("D51").Value = "0.0002007"
$ws.Range("E51").Value = "0.26%"
$ws.Range("G51").Value = "6"

# Updated coin names and links (plain text, no numeric coercion risk)
$ws.Range("B18").Value = "HotbitToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("B19").Value = "LEO"
$ws.Range("C19").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("B21").Value = "MCDex"
$ws.Range("C21").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("B22").Value = "ProBitToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("B23").Value = "ZBToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("B24").Value = "CoinExToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("B25").Value = "BitKan"
$ws.Range("C25").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
